# Insert a new weekly price record for "Apio" (Macroferia Regional de Talca)
# at row 164 of Sheet1, pushing the existing rows 164:293 down to 165:294.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 164 through the end down by one row, creating an empty row 164.
$ws.Rows("164:164").Insert()

# Populate the new row 164 with the new record's data.
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 45062
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = 100112017
$ws.Range("G164").Value = "Apio"
$ws.Range("H164").Value = "Americana (o)"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 700
$ws.Range("K164").Value = 6000
$ws.Range("L164").Value = 6000
$ws.Range("M164").Value = 6000
$ws.Range("N164").Value = "`$/docena de matas"
$ws.Range("O164").Value = "Provincia del Elquí"
$ws.Range("P164").Value = 1000
$ws.Range("Q164").Value = 6
$ws.Range("R164").Value = "Hortaliza"
